$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: a second property-building entry ("#5: property building done"),
# laid out like row 1 but shifted one column left with a leading count
# cell in A2.
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "台新國際商業銀行南京東路分行"
$ws.Range("C2").Value = "活期存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "李慶華"
$ws.Range("F2").Value = 7175

# A2 picks up the same bordered/centered style already used by row 1
# (copy the fully-resolved format in one shot so it collapses onto the
# existing style record instead of minting a new one).
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# B2:F2 get a distinct (but visually default) style record of their own,
# separate from row 1's bordered style and from the sheet's base style.
$ws.Range("B2:F2").NumberFormat = "General"
